$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the curated motif table (A3:G16) by column A ascending. It had
# previously been sorted by column F (position) descending; the new sort
# key is the added numbering/order column A.
$sortRange = $ws.Range("A3:G16")
$keyRange  = $ws.Range("A3:A16")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# The engine's sort carries an empty styled band cell into F12 (a position
# that never had an F value for that particular data row). Drop it so the
# row matches the source record exactly (no F cell at all).
$ws.Range("F12").Clear()

# Leave the selection where the user ended up after making the edit.
$ws.Range("B16").Select()
